$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "244.46"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "23.09"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.419"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.05952"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.388"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.9252"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1416"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07437"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.03415"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03055"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.09342"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.938"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001595"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.04806"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0005945"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.005707"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.004152"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0009845"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.00007706"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.659"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.3242"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1341"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.03923"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006215"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1074"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002612"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007322"
$ws.Range("E44").Value = "43LocalTradersLCTBestin24h"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00005165"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0005805"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
